$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("destination_simple")

# --- Update the B-column category labels (re-mapping the old "Other"/
# --- "Social/Recreation" buckets into new, more granular categories) ---
$ws.Range("B4").Value  = "Drop off/Pick up"
$ws.Range("B10").Value = "Errands"
$ws.Range("B11").Value = "Health and Exercise"
$ws.Range("B12").Value = "Social/Recreation/Eat Meal"
$ws.Range("B13").Value = "Health and Exercise"
$ws.Range("B14").Value = "Social/Recreation/Eat Meal"
$ws.Range("B15").Value = "Social/Recreation/Eat Meal"
$ws.Range("B16").Value = "Social/Recreation/Eat Meal"
$ws.Range("B17").Value = "Social/Recreation/Eat Meal"
$ws.Range("B19").Value = "Errands"
$ws.Range("B20").Value = "Social/Recreation/Eat Meal"

# --- Header cell: the stale "Social/Recreation" shared string is dropped from
# --- the table entirely, so the header simply re-asserts "Destination purpose" ---
$ws.Range("A1").Value = "Destination purpose"

# --- Widen column A to fit the longer category text, and leave the cursor on
# --- the cell that was being edited ---
$ws.Columns.Item(1).ColumnWidth = 49
$ws.Range("F8").Select()
